$d = $word.ActiveDocument
$table = $d.Tables.Item(6)
$cell = $table.Cell(1, 2)
$cell.Range.Text = "2023-04-12"

$table2 = $d.Tables.Item(6)
$cell2 = $table2.Cell(1, 2)
$cell2.Range.Font.Name = "Calibri"

$table3 = $d.Tables.Item(6)
$cell3 = $table3.Cell(1, 2)
$cell3.Range.Font.NameAscii = "Calibri"
